$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell B11 on the "Rules" sheet currently holds the text "R40".
# It is updated to hold the text "1" (kept as text, not converted to a
# number, mirroring the rest of column B which stores rule codes as
# strings). Force text storage via NumberFormat "@" so the literal "1"
# isn't auto-coerced into a numeric value.
$cell = $ws.Range("B11")
$cell.NumberFormat = "@"
$cell.Value = "1"
